$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.831.42'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '1.721.64'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('D4').Value = "'0.9987"
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = "'239.68"
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('D6').Value = "'0.9995"
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').Value = "'0.4759"
$ws.Range('E7').Value = '  -2.17%  '
$ws.Range('D8').Value = "'0.2559"
$ws.Range('E8').Value = '  -0.93%  '
$ws.Range('D9').Value = "'0.06115"
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').Value = '1.720.05'
$ws.Range('E10').Value = '  +0.32%  '
$ws.Range('D11').Value = "'15.84"
$ws.Range('E11').Value = '  +2.21%  '
$ws.Range('D12').Value = "'0.06889"
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('D13').Value = "'0.5959"
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').Value = "'4.396"
$ws.Range('E14').Value = '  -1.64%  '
$ws.Range('D15').Value = "'76.32"
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').Value = "'1.000"
$ws.Range('D17').Value = '26.724.29'
$ws.Range('E17').Value = '  +1.82%  '
$ws.Range('D18').Value = "'0.9984"
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = "'0.000006981"
$ws.Range('E19').Value = '  -1.62%  '
$ws.Range('D20').Value = "'11.26"
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').Value = '1.938.19'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = "'4.362"
$ws.Range('E22').Value = '  -1.00%  '
$ws.Range('D23').Value = "'8.325"
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('D24').Value = "'5.050"
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = "'140.77"
$ws.Range('E25').Value = '  +3.24%  '
$ws.Range('D26').Value = "'15.14"
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = "'1.792"
$ws.Range('E27').Value = '  +3.57%  '
$ws.Range('D28').Value = "'106.12"
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('D29').Value = "'1.375"
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('D30').Value = "'3.935"
$ws.Range('E30').Value = '  +1.56%  '
$ws.Range('D31').Value = "'0.07887"
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('D32').Value = "'3.628"
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('D33').Value = "'0.04618"
$ws.Range('E33').Value = '  +4.27%  '
$ws.Range('D34').Value = "'2.595"
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').Value = "'0.9936"
$ws.Range('E35').Value = '  +0.30%  '
$ws.Range('D36').Value = "'0.6093"
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('D37').Value = "'0.9154"
$ws.Range('E37').Value = '  -2.17%  '
$ws.Range('D38').Value = "'2.506"
$ws.Range('E38').Value = '  +5.58%  '
$ws.Range('D39').Value = "'1.972"
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('D40').Value = "'0.9988"
$ws.Range('E40').Value = '  +0.31%  '
$ws.Range('D41').Value = "'5.640"
$ws.Range('E41').Value = '  +4.95%  '
$ws.Range('D42').Value = "'0.01480"
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('D43').Value = "'99.50"
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('D44').Value = "'0.3779"
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('D45').Value = "'6.701"
$ws.Range('E45').Value = '  -1.74%  '
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('D47').Value = "'0.05338"
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').Value = "'7.691"
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('D49').Value = "'29.71"
$ws.Range('E49').Value = '  -2.82%  '
$ws.Range('D50').Value = "'1.230"
$ws.Range('E50').Value = '  +1.65%  '
$ws.Range('D51').Value = "'1.002"
$ws.Range('E51').Value = '  +0.27%  '
